$wb = $excel.ActiveWorkbook

# --- Sheet "29.12.23" (index 2): add a new paid-off record + a partial-transfer note ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A17").Value2 = "Đã trả"
$ws2.Range("B17").Value2 = 70000000
$ws2.Range("C17").Value2 = 1
$ws2.Range("D17").Formula = "=B17*C17"
$ws2.Range("E17").Value2 = "Chuyển ngày 5.2.24"

$ws2.Range("D18").Formula = "=SUM(D14:D17)"

$ws2.Range("E21").Value2 = "Chuyển sang công nợ ngày 6.1.24"
$ws2.Range("G21").NumberFormat = "#,##0"

$ws2.Range("D21").Select()

# --- Sheet "6.1.24" (index 3): add the matching paid-back record ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A14").Value2 = "Đã trả ngày 5.2.24"
$ws3.Range("B14").Value2 = 6750000
$ws3.Range("C14").Value2 = 1
$ws3.Range("D14").Formula = "=B14*C14"

$ws3.Range("D15").Formula = "=SUM(D11:D14)"

$ws3.Activate()
$ws3.Range("D16").Select()
